$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Alyssa 4, Lauren 8, Kayla 5, Madison 14"
$ws.Range("E3").Value = "Lauren 12, Lauren 11"
$ws.Range("E4").Value = "Lauren 8, Alyssa 4"
$ws.Range("E5").Value = "Elizabeth 10, Emma 9"
$ws.Range("E6").Value = "Elizabeth 10, Emma 9"
$ws.Range("E7").Value = "Elizabeth 10, Emma 9"
$ws.Range("E8").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("E11").Value = "Alyssa 4, Kayla 5, Madison 14"
$ws.Range("E12").Value = "Samantha 18"
$ws.Range("E13").Value = "Olivia 1"
$ws.Range("E14").Value = "Kayla 5, Isabella 3, Madison 14, Taylor 17, Alyssa 4, Samantha 18"
$ws.Range("E15").Value = "Sarah 2, Kayla 15, Lauren 12"
$ws.Range("E16").Value = "Alyssa 4"
$ws.Range("E17").Value = "Lauren 8, Kayla 15"
$ws.Range("E18").Value = "Isabella 3"
$ws.Range("E19").Value = "Lauren 8, Lauren 11"
$ws.Range("E20").ClearContents()
$ws.Range("E21").Value = "Kayla 13"
$ws.Range("E22").ClearContents()
$ws.Range("E23").Value = "Hannah 16"
$ws.Range("E24").Value = "Kayla 7, Kayla 13"
$ws.Range("E25").Value = "Kayla 15"
$ws.Range("E26").Value = "Kayla 7"
$ws.Range("E27").Value = "Emma 9"
$ws.Range("E28").Value = "Kayla 15, Kayla 5"
$ws.Range("E29").Value = "Kayla 15, Taylor 17"
$ws.Range("E30").Value = "Taylor 6, Elizabeth 10"
$ws.Range("E31").Value = "Hannah 16, Isabella 3"
$ws.Range("E32").ClearContents()
$ws.Range("E33").ClearContents()
$ws.Range("E34").Value = "Kayla 5, Emma 9"
$ws.Range("E35").Value = "Alyssa 4, Lauren 11"
$ws.Range("E36").ClearContents()
$ws.Range("E37").Value = "Sarah 2"
$ws.Range("E38").ClearContents()
$ws.Range("E40").Value = "Kayla 5"
$ws.Range("E41").Value = "Elizabeth 10, Alyssa 4"
$ws.Range("E42").Value = "Lauren 12, Lauren 11, Emma 9, Sarah 2, Olivia 1"
$ws.Range("E43").ClearContents()
$ws.Range("E44").ClearContents()
$ws.Range("E45").Value = "Hannah 16, Olivia 1"
$ws.Range("E46").Value = "Taylor 6"
$ws.Range("E47").Value = "Taylor 6"
$ws.Range("E48").Value = "Taylor 6"
$ws.Range("E49").ClearContents()
$ws.Range("E50").Value = "Emma 9"
$ws.Range("E51").Value = "Lauren 12, Sarah 2, Hannah 16"
$ws.Range("E52").Value = "Lauren 11"
$ws.Range("E53").Value = "Lauren 8, Kayla 13, Isabella 3, Kayla 7"
$ws.Range("E55").Value = "Madison 14, Taylor 17, Kayla 5"
$ws.Range("E57").ClearContents()
$ws.Range("E58").Value = "Sarah 2"
$ws.Range("E59").Value = "Olivia 1"
$ws.Range("E60").Value = "Olivia 1"
$ws.Range("E61").Value = "Olivia 1"
$ws.Range("E62").Value = "Samantha 18"
$ws.Range("E63").Value = "Samantha 18"
$ws.Range("E64").Value = "Samantha 18"
$ws.Range("E65").Value = "Alyssa 4, Kayla 15, Hannah 16"
$ws.Range("E66").ClearContents()
$ws.Range("E67").Value = "Hannah 16, Samantha 18"
$ws.Range("E68").ClearContents()
$ws.Range("E69").Value = "Kayla 7, Taylor 17"
$ws.Range("E70").Value = "Kayla 7, Taylor 17"
$ws.Range("E71").Value = "Lauren 12, Taylor 6"
$ws.Range("E72").Value = "Sarah 2, Isabella 3, Elizabeth 10, Madison 14"
$ws.Range("E73").Value = "Kayla 13"
$ws.Range("E74").ClearContents()
$ws.Range("E75").ClearContents()
$ws.Range("E76").Value = "Lauren 12"
$ws.Range("E77").Value = "Isabella 3"
$ws.Range("E78").Value = "Lauren 12, Kayla 15, Lauren 11, Elizabeth 10"
$ws.Range("E79").ClearContents()
$ws.Range("E80").Value = "Kayla 13, Taylor 17"
$ws.Range("E81").Value = "Lauren 12, Kayla 15, Lauren 11, Elizabeth 10, Isabella 3"
$ws.Range("E82").ClearContents()
$ws.Range("E83").ClearContents()
$ws.Range("E84").Value = "Kayla 13, Taylor 6, Madison 14"
$ws.Range("E85").Value = "Kayla 13, Lauren 11, Kayla 7"
$ws.Range("E86").Value = "Lauren 8, Taylor 17"
$ws.Range("E87").Value = "Sarah 2, Madison 14"
$ws.Range("E88").Value = "Kayla 5"
$ws.Range("E89").ClearContents()
$ws.Range("E90").Value = "Isabella 3, Hannah 16"
$ws.Range("E91").ClearContents()
$ws.Range("E92").Value = "Taylor 6, Kayla 13"
$ws.Range("E93").Value = "Lauren 8, Sarah 2"
$ws.Range("E94").Value = "Olivia 1, Emma 9, Hannah 16, Kayla 7"
$ws.Range("E95").Value = "Lauren 8, Taylor 6, Samantha 18, Madison 14, Taylor 17"
$ws.Range("E96").Value = "Olivia 1, Samantha 18, Kayla 7"
$ws.Range("E97").ClearContents()
